# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner market-data refresh to Siren_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 144.46666  # H11
$ws.Cells.Item(11, 9).Value = 144.46666  # I11
$ws.Cells.Item(11, 11).Value = 144.46666  # K11
$ws.Cells.Item(11, 13).Value = -4.46665999999999  # M11
$ws.Cells.Item(29, 8).Value = 903  # H29
$ws.Cells.Item(29, 10).Value = 0  # J29
$ws.Cells.Item(29, 12).Value = 0  # L29
$ws.Cells.Item(29, 14).ClearContents()  # N29
$ws.Cells.Item(32, 8).Value = 5547.476  # H32
$ws.Cells.Item(32, 10).Value = 5377.1816  # J32
$ws.Cells.Item(32, 12).Value = 5377.1816  # L32
$ws.Cells.Item(32, 14).Value = -6029.1816  # N32
$ws.Cells.Item(41, 8).Value = 149.5  # H41
$ws.Cells.Item(41, 10).Value = 100  # J41
$ws.Cells.Item(41, 12).Value = 100  # L41
$ws.Cells.Item(41, 14).Value = -980  # N41
$ws.Cells.Item(64, 8).Value = 45734.3  # H64
$ws.Cells.Item(64, 9).Value = 71224.664  # I64
$ws.Cells.Item(64, 10).Value = 7498.75  # J64
$ws.Cells.Item(64, 11).Value = 71224.664  # K64
$ws.Cells.Item(64, 12).Value = 7498.75  # L64
$ws.Cells.Item(64, 13).Value = -70976.664  # M64
$ws.Cells.Item(64, 14).Value = -7994.75  # N64
$ws.Cells.Item(67, 8).Value = 45734.3  # H67
$ws.Cells.Item(67, 9).Value = 71224.664  # I67
$ws.Cells.Item(67, 10).Value = 7498.75  # J67
$ws.Cells.Item(67, 11).Value = 71224.664  # K67
$ws.Cells.Item(67, 12).Value = 7498.75  # L67
$ws.Cells.Item(67, 13).Value = -70366.664  # M67
$ws.Cells.Item(67, 14).Value = -9214.75  # N67
$ws.Cells.Item(97, 8).Value = 1878.8823  # H97
$ws.Cells.Item(97, 9).Value = 1500  # I97
$ws.Cells.Item(97, 10).Value = 2420.1428  # J97
$ws.Cells.Item(97, 11).Value = 4500  # K97
$ws.Cells.Item(97, 12).Value = 7260.428400000001  # L97
$ws.Cells.Item(97, 13).Value = -4004  # M97
$ws.Cells.Item(97, 14).Value = -8252.428400000001  # N97
$ws.Cells.Item(116, 8).Value = 859288.9399999999  # H116
$ws.Cells.Item(116, 9).Value = 1854809.9  # I116
$ws.Cells.Item(116, 11).Value = 1854809.9  # K116
$ws.Cells.Item(116, 13).Value = -1851367.9  # M116
$ws.Cells.Item(125, 8).Value = 3899.2  # H125
$ws.Cells.Item(125, 10).Value = 3899.2  # J125
$ws.Cells.Item(125, 12).Value = 35092.8  # L125
$ws.Cells.Item(125, 14).Value = -40012.8  # N125
$ws.Cells.Item(135, 8).Value = 9544.75  # H135
$ws.Cells.Item(135, 9).Value = 12342.75  # I135
$ws.Cells.Item(135, 11).Value = 111084.75  # K135
$ws.Cells.Item(135, 13).Value = -108549.75  # M135
$ws.Cells.Item(138, 8).Value = 3268.125  # H138
$ws.Cells.Item(138, 9).Value = 698.4  # I138
$ws.Cells.Item(138, 10).Value = 4436.1816  # J138
$ws.Cells.Item(138, 11).Value = 2095.2  # K138
$ws.Cells.Item(138, 12).Value = 13308.5448  # L138
$ws.Cells.Item(138, 13).Value = 3044.8  # M138
$ws.Cells.Item(138, 14).Value = -23588.5448  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5907.6294  # H61
$ws.Cells.Item(61, 9).Value = 6716.727  # I61
$ws.Cells.Item(61, 10).Value = 5351.375  # J61
$ws.Cells.Item(61, 11).Value = 6716.727  # K61
$ws.Cells.Item(61, 12).Value = 5351.375  # L61
$ws.Cells.Item(61, 13).Value = -6504.727  # M61
$ws.Cells.Item(61, 14).Value = -5775.375  # N61
$ws.Cells.Item(102, 8).Value = 6111  # H102
$ws.Cells.Item(102, 9).Value = 6062.375  # I102
$ws.Cells.Item(102, 11).Value = 6062.375  # K102
$ws.Cells.Item(102, 13).Value = -4440.375  # M102
$ws.Cells.Item(122, 8).Value = 566183.75  # H122
$ws.Cells.Item(122, 9).Value = 5230.6924  # I122
$ws.Cells.Item(122, 11).Value = 15692.0772  # K122
$ws.Cells.Item(122, 13).Value = -13242.0772  # M122
$ws.Cells.Item(136, 8).Value = 5907.6294  # H136
$ws.Cells.Item(136, 9).Value = 6716.727  # I136
$ws.Cells.Item(136, 10).Value = 5351.375  # J136
$ws.Cells.Item(136, 11).Value = 20150.181  # K136
$ws.Cells.Item(136, 12).Value = 16054.125  # L136
$ws.Cells.Item(136, 13).Value = -17600.181  # M136
$ws.Cells.Item(136, 14).Value = -21154.125  # N136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 13848.267  # H99
$ws.Cells.Item(99, 9).Value = 14399.929  # I99
$ws.Cells.Item(99, 11).Value = 14399.929  # K99
$ws.Cells.Item(99, 13).Value = -12901.929  # M99
$ws.Cells.Item(132, 8).Value = 53796  # H132
$ws.Cells.Item(132, 10).Value = 53796  # J132
$ws.Cells.Item(132, 12).Value = 53796  # L132
$ws.Cells.Item(132, 14).Value = -63916  # N132

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2833.3333  # H16
$ws.Cells.Item(16, 9).Value = 2833.3333  # I16
$ws.Cells.Item(16, 11).Value = 2833.3333  # K16
$ws.Cells.Item(16, 13).Value = -2546.3333  # M16
$ws.Cells.Item(31, 8).Value = 4674.5  # H31
$ws.Cells.Item(31, 9).Value = 1349  # I31
$ws.Cells.Item(31, 11).Value = 1349  # K31
$ws.Cells.Item(31, 13).Value = -1054  # M31
$ws.Cells.Item(34, 8).Value = 4674.5  # H34
$ws.Cells.Item(34, 9).Value = 1349  # I34
$ws.Cells.Item(34, 11).Value = 1349  # K34
$ws.Cells.Item(34, 13).Value = -1147  # M34
$ws.Cells.Item(62, 8).Value = 82118  # H62
$ws.Cells.Item(62, 10).Value = 153523.75  # J62
$ws.Cells.Item(62, 12).Value = 153523.75  # L62
$ws.Cells.Item(62, 14).Value = -154771.75  # N62
$ws.Cells.Item(65, 8).Value = 82118  # H65
$ws.Cells.Item(65, 10).Value = 153523.75  # J65
$ws.Cells.Item(65, 12).Value = 767618.75  # L65
$ws.Cells.Item(65, 14).Value = -773858.75  # N65
$ws.Cells.Item(113, 8).Value = 2833.3333  # H113
$ws.Cells.Item(113, 9).Value = 2833.3333  # I113
$ws.Cells.Item(113, 11).Value = 2833.3333  # K113
$ws.Cells.Item(113, 13).Value = -663.3332999999998  # M113

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 572.2222  # H23
$ws.Cells.Item(23, 9).Value = 478  # I23
$ws.Cells.Item(23, 11).Value = 1434  # K23
$ws.Cells.Item(23, 13).Value = -1199  # M23
$ws.Cells.Item(51, 8).Value = 1808.2307  # H51
$ws.Cells.Item(51, 9).Value = 995.625  # I51
$ws.Cells.Item(51, 10).Value = 3108.4  # J51
$ws.Cells.Item(51, 11).Value = 2986.875  # K51
$ws.Cells.Item(51, 12).Value = 9325.200000000001  # L51
$ws.Cells.Item(51, 13).Value = -2526.875  # M51
$ws.Cells.Item(51, 14).Value = -10245.2  # N51
$ws.Cells.Item(117, 8).Value = 4678.4  # H117
$ws.Cells.Item(117, 10).Value = 4700  # J117
$ws.Cells.Item(117, 12).Value = 14100  # L117
$ws.Cells.Item(117, 14).Value = -20984  # N117
$ws.Cells.Item(129, 8).Value = 2072.1667  # H129
$ws.Cells.Item(129, 9).Value = 833.3333  # I129
$ws.Cells.Item(129, 10).Value = 3311  # J129
$ws.Cells.Item(129, 11).Value = 2499.9999  # K129
$ws.Cells.Item(129, 12).Value = 9933  # L129
$ws.Cells.Item(129, 13).Value = 2500.0001  # M129
$ws.Cells.Item(129, 14).Value = -19933  # N129
$ws.Cells.Item(131, 8).Value = 6632.2144  # H131
$ws.Cells.Item(131, 9).Value = 9981.625  # I131
$ws.Cells.Item(131, 10).Value = 2166.3333  # J131
$ws.Cells.Item(131, 11).Value = 29944.875  # K131
$ws.Cells.Item(131, 12).Value = 6498.999899999999  # L131
$ws.Cells.Item(131, 13).Value = -24904.875  # M131
$ws.Cells.Item(131, 14).Value = -16578.9999  # N131
$ws.Cells.Item(136, 8).Value = 976.6667  # H136
$ws.Cells.Item(136, 9).Value = 976.6667  # I136
$ws.Cells.Item(136, 11).Value = 2930.0001  # K136
$ws.Cells.Item(136, 13).Value = 2169.9999  # M136
$ws.Cells.Item(139, 8).Value = 1430245.9  # H139
$ws.Cells.Item(139, 9).Value = 1580271.8  # I139
$ws.Cells.Item(139, 11).Value = 4740815.4  # K139
$ws.Cells.Item(139, 13).Value = -4735675.4  # M139
$ws.Cells.Item(140, 8).Value = 14769.429  # H140
$ws.Cells.Item(140, 9).Value = 17843  # I140
$ws.Cells.Item(140, 10).Value = 3499.6667  # J140
$ws.Cells.Item(140, 11).Value = 53529  # K140
$ws.Cells.Item(140, 12).Value = 10499.0001  # L140
$ws.Cells.Item(140, 13).Value = -48349  # M140
$ws.Cells.Item(140, 14).Value = -20859.0001  # N140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(75, 8).Value = 35000  # H75
$ws.Cells.Item(75, 10).Value = 35000  # J75
$ws.Cells.Item(75, 12).Value = 35000  # L75
$ws.Cells.Item(75, 14).Value = -36748  # N75
$ws.Cells.Item(78, 8).Value = 35000  # H78
$ws.Cells.Item(78, 10).Value = 35000  # J78
$ws.Cells.Item(78, 12).Value = 105000  # L78
$ws.Cells.Item(78, 14).Value = -113736  # N78
$ws.Cells.Item(126, 8).Value = 14821.875  # H126
$ws.Cells.Item(126, 9).Value = 33859.6  # I126
$ws.Cells.Item(126, 11).Value = 101578.8  # K126
$ws.Cells.Item(126, 13).Value = -99108.79999999999  # M126

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 31991.625  # H40
$ws.Cells.Item(40, 9).Value = 53574.875  # I40
$ws.Cells.Item(40, 10).Value = 10408.375  # J40
$ws.Cells.Item(40, 11).Value = 53574.875  # K40
$ws.Cells.Item(40, 12).Value = 10408.375  # L40
$ws.Cells.Item(40, 13).Value = -53438.875  # M40
$ws.Cells.Item(40, 14).Value = -10680.375  # N40
$ws.Cells.Item(61, 8).Value = 6510.7617  # H61
$ws.Cells.Item(61, 9).Value = 4426.6875  # I61
$ws.Cells.Item(61, 10).Value = 13179.8  # J61
$ws.Cells.Item(61, 11).Value = 4426.6875  # K61
$ws.Cells.Item(61, 12).Value = 13179.8  # L61
$ws.Cells.Item(61, 13).Value = -4224.6875  # M61
$ws.Cells.Item(61, 14).Value = -13583.8  # N61
$ws.Cells.Item(68, 8).Value = 5100  # H68
$ws.Cells.Item(71, 8).Value = 5100  # H71
$ws.Cells.Item(113, 8).Value = 6510.7617  # H113
$ws.Cells.Item(113, 9).Value = 4426.6875  # I113
$ws.Cells.Item(113, 10).Value = 13179.8  # J113
$ws.Cells.Item(113, 11).Value = 4426.6875  # K113
$ws.Cells.Item(113, 12).Value = 13179.8  # L113
$ws.Cells.Item(113, 13).Value = -2256.6875  # M113
$ws.Cells.Item(113, 14).Value = -17519.8  # N113
$ws.Cells.Item(122, 8).Value = 3925.5186  # H122
$ws.Cells.Item(122, 9).Value = 3503.8235  # I122
$ws.Cells.Item(122, 10).Value = 4642.4  # J122
$ws.Cells.Item(122, 11).Value = 10511.4705  # K122
$ws.Cells.Item(122, 12).Value = 13927.2  # L122
$ws.Cells.Item(122, 13).Value = -8061.470499999999  # M122
$ws.Cells.Item(122, 14).Value = -18827.2  # N122

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 340812.75  # H62
$ws.Cells.Item(62, 9).Value = 678625.5  # I62
$ws.Cells.Item(62, 11).Value = 678625.5  # K62
$ws.Cells.Item(62, 13).Value = -678001.5  # M62
$ws.Cells.Item(65, 8).Value = 340812.75  # H65
$ws.Cells.Item(65, 9).Value = 678625.5  # I65
$ws.Cells.Item(65, 11).Value = 3393127.5  # K65
$ws.Cells.Item(65, 13).Value = -3390007.5  # M65

